# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Swap the province names between row 47 (Lugo) and row 48 (Almeria)
$ws.Range("A47").Value = "Almeria"
$ws.Range("A48").Value = "Lugo"

# Swap the "Casos activos" values between row 47 and row 48
$ws.Range("C47").Value = 72
$ws.Range("C48").Value = 5

# Update the "last updated" timestamp text
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 05:46"
